$wb = $excel.ActiveWorkbook

$wsSheet1 = $wb.Worksheets.Item('Sheet 1')
$wsReferences = $wb.Worksheets.Item('References')

$wsSheet1.Range('E5').Value = 'Based on official disease reports to the WOAH'
$wsSheet1.Range('E6').Value = 'CCPP is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005.'
$wsSheet1.Range('E7').Value = 'As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:'
$wsSheet1.Range('E14').Value = 'Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently.'
$wsSheet1.Range('E17').Value = 'For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}.'
$wsSheet1.Range('E21').Value = 'A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}.'
$wsSheet1.Range('E34').Value = 'Humans are not susceptible,  and therefore there is no direct impact on public health ({ref008:WOAH Technical disease card}).'
$wsSheet1.Range('E42').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the bacterium characteristics. '
$wsSheet1.Range('E56').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters.'
$wsSheet1.Range('E68').Value = 'WOAH-prescribed tests for international trade include:the commercial competitive enzyme-linked immunosorbent assay (ELISA) and in-house indirect ELISA ({ref010:WOAH, Terrestrial Manual})'
$wsSheet1.Range('E93').Value = 'Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data.'
$wsReferences.Range('C2').Value = 'WOAH-WAHIS (WOAH World Animal Health Information System)'
$wsReferences.Range('C6').Value = 'WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$wsReferences.Range('C9').Value = 'WOAH (World Organisation for Animal Health) Technical Disease Card: Contagious caprine pleuropneumonia, 2009'
$wsReferences.Range('C10').Value = 'WOAH (World Organisation for Animal Health), 2021. Contagious caprine pleuropneumonia. Chapter 14.3. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$wsReferences.Range('C11').Value = 'WOAH (World Organisation for Animal Health), 2021. Contagious caprine pleuropneumonia. Chapter 3.08.04. WOAH Terrestrial Manual 2021. WOAH, Paris, France'
